$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the description text for the "BRL" row (row 15, column C) to the
# finalized wording used in the data release README.
$ws.Range("C15").Value = "Whether the concentration is below the reporting limit. Yes indicates the concentration is above the method detection limit, but below the reporting limit"

# Reflect the cursor/selection position left behind after the edit.
$ws.Range("C16").Select()
